# Apply cryptos list price/volume update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.716.53'
$ws.Range("E2").Value = '  -2.39%  '
$ws.Range("D3").Value = '2.292.80'
$ws.Range("E3").Value = '  -5.24%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'546.93"
$ws.Range("E5").Value = '  -1.33%  '
$ws.Range("D6").Value = "'130.82"
$ws.Range("E6").Value = '  -4.79%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = "'0.573"
$ws.Range("E8").Value = '  -2.71%  '
$ws.Range("D9").Value = '2.291.41'
$ws.Range("E9").Value = '  -5.30%  '
$ws.Range("E10").Value = '  -3.52%  '
$ws.Range("D11").Value = "'5.55"
$ws.Range("E11").Value = '  -2.84%  '
$ws.Range("E12").Value = '  +0.87%  '
$ws.Range("D13").Value = "'0.335"
$ws.Range("E13").Value = '  -5.12%  '
$ws.Range("D14").Value = "'23.82"
$ws.Range("E14").Value = '  -4.42%  '
$ws.Range("D15").Value = '2.701.47'
$ws.Range("E15").Value = '  -5.35%  '
$ws.Range("D16").Value = '58.665.46'
$ws.Range("E16").Value = '  -2.35%  '
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("E17").Value = '  -3.34%  '
$ws.Range("D18").Value = '2.276.20'
$ws.Range("E18").Value = '  -5.71%  '
$ws.Range("D20").Value = "'4.30"
$ws.Range("E20").Value = '  -4.60%  '
$ws.Range("D21").Value = "'315.29"
$ws.Range("E21").Value = '  -3.75%  '
$ws.Range("E22").Value = '  -4.23%  '
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = '  -0.20%  '
$ws.Range("D24").Value = "'63.03"
$ws.Range("E24").Value = '  -3.21%  '
$ws.Range("D25").Value = "'0.173"
$ws.Range("E25").Value = '  -2.64%  '
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").Value = "'8.09"
$ws.Range("E27").Value = '  -7.25%  '
$ws.Range("D28").Value = "'1.31"
$ws.Range("E28").Value = '  -5.37%  '
$ws.Range("D29").Value = "'1.74"
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("D30").Value = "'169.83"
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("D31").Value = '0.0₃0728'
$ws.Range("E31").Value = '  -6.12%  '
$ws.Range("E32").Value = '  -5.15%  '
$ws.Range("D33").Value = "'1.07"
$ws.Range("E33").Value = '  +0.44%  '
$ws.Range("D34").Value = "'0.384"
$ws.Range("E34").Value = '  -5.02%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = "'17.82"
$ws.Range("E36").Value = '  -3.97%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("E38").Value = '  -6.04%  '
$ws.Range("D39").Value = "'3.96"
$ws.Range("E39").Value = '  -6.20%  '
$ws.Range("E40").Value = '  -2.24%  '
$ws.Range("D41").Value = "'1.51"
$ws.Range("E41").Value = '  -5.45%  '
$ws.Range("D42").Value = "'299.46"
$ws.Range("E42").Value = '  -8.45%  '
$ws.Range("D43").Value = "'140.44"
$ws.Range("E43").Value = '  -2.95%  '
$ws.Range("D44").Value = "'3.45"
$ws.Range("E44").Value = '  -5.86%  '
$ws.Range("D45").Value = "'0.0954"
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("D46").Value = "'0.0501"
$ws.Range("E46").Value = '  -3.32%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = "'0.557"
$ws.Range("E47").Value = '  -3.45%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").Value = "'18.59"
$ws.Range("E48").Value = '  -7.33%  '
$ws.Range("D49").Value = "'0.0215"
$ws.Range("E49").Value = '  -3.63%  '
$ws.Range("D50").Value = "'16.66"
$ws.Range("E50").Value = '  -4.87%  '
$ws.Range("D51").Value = "'11.00"
$ws.Range("E51").Value = '  -0.31%  '
